# ----------------------------------------------------------------------
# Applies the tracked edits to ExcelHW_LisaCannon.docx
# ----------------------------------------------------------------------
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Given the provided data ... " question: merge the two runs
#    that used to be split by the _GoBack bookmark into a single run
#    (the bookmark itself moves into the new default header, see below).
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Given the provided data, what are three conclusions we can draw about Kickstarter campaigns?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Given the provided data, what are three conclusions we can draw about Kickstarter campaigns?", 2)

# ------------------------------------------------------------------
# 2) "... highest proportion of successes is music" -> "... are music"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "highest proportion of successes is music, theater, and film/video.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "highest proportion of successes are music, theater, and film/video.", 2)

# ------------------------------------------------------------------
# 3) "over 300,000 projects which have raised" -> "over 300,000 projects in Kickstarter which have raised"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "over 300,000 projects which have raised over $2 billion",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "over 300,000 projects in Kickstarter which have raised over $2 billion", 2)

# ------------------------------------------------------------------
# 4) Append a new sentence after "... than is present in the population. "
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "than is present in the population. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found = $d.Content
$found.Find.Execute(
    "than is present in the population. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($found.End, $found.End)
$insertPoint.InsertAfter("Perhaps changing the sampling technique could ensure random selection of projects. ")

# ------------------------------------------------------------------
# 5) "1/3 of the average number of backers" -> "1/3 of the mean number of backers"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "1/3 of the average number of backers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1/3 of the mean number of backers", 2)

Write-Host "Body text edits applied"
